$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("decomposition_main_te")

# Update the three header/label cells whose text actually changed.
$ws.Range("A6").Value = "Mandatory structured"
$ws.Range("A8").Value = "Choice "
$ws.Range("E3").Value = 'Def$\times$Ppl pymnt'

# Column A needs to widen to fit the new, longer label ("Mandatory structured").
$ws.Columns.Item(1).ColumnWidth = 18.6
